$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 550, shifting existing rows 550-611 down to 551-612.
$ws.Rows.Item(550).Insert()

# Populate the newly inserted row 550 with the new data record.
$ws.Cells.Item(550, 1).Value = 11
$ws.Cells.Item(550, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(550, 3).Value = "Bíobío"
$ws.Cells.Item(550, 4).Value = 45132
$ws.Cells.Item(550, 5).Value = 8
$ws.Cells.Item(550, 6).Value = 100112006
$ws.Cells.Item(550, 7).Value = "Repollo"
$ws.Cells.Item(550, 8).Value = "Crespo record"
$ws.Cells.Item(550, 9).Value = "Primera"
$ws.Cells.Item(550, 10).Value = 1000
$ws.Cells.Item(550, 11).Value = 900
$ws.Cells.Item(550, 12).Value = 1000
$ws.Cells.Item(550, 13).Value = 950
$ws.Cells.Item(550, 14).Value = "$/unidad"
$ws.Cells.Item(550, 15).Value = "Región Metropolitana"
$ws.Cells.Item(550, 16).Value = 950
$ws.Cells.Item(550, 17).Value = 1
$ws.Cells.Item(550, 18).Value = "Hortaliza"
